$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    0.03156792118756969,
    0.939990559993592,
    0.05167416594198131,
    0.04147550293697191,
    0.03263783550394532,
    0.04514868199427455,
    0.07437797941466077,
    0.9467563311337541,
    0.058308350443767,
    0.9376882015336185,
    0.06376744125821353,
    0.06029205332833725,
    0.05779998098317343,
    0.05616426305477923,
    0.05472168419632627,
    0.955677609879598,
    0.9355315004666163,
    0.04800921162175867,
    0.04565151322148366,
    0.9394860282568145,
    0.05064639480204739,
    0.04166514918221484,
    0.04501212832098334,
    0.04504873063684611,
    0.03571889043675101,
    0.9420327884994528,
    0.04687582373724237,
    0.05125038313061273,
    0.05632607957537897,
    0.9463047396100521,
    0.0563236159452473,
    0.06160953389031205,
    0.05356213928576739,
    0.06610446550074273,
    0.06835630603067969,
    0.05638666982658216,
    0.0587111191203309,
    0.06621844033311444,
    0.9563735210504116,
    0.06277490953660118
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $values[$i]
}

$wb.Save()
